$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "A1: $($ws.Cells.Item(1,1).Value())"
Write-Host "E1: $($ws.Cells.Item(1,5).Value())"
Write-Host "F1: $($ws.Cells.Item(1,6).Value())"
